$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTICE_DATE / REPORT_DATE are stored as plain text strings in this sheet
$ws.Range("M2").Value = "2020-12-24 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = -44125029.51
$ws.Range("P2").Value = 143337667.45
$ws.Range("Q2").Value = 180305197.21

$ws.Range("S2").Value = 114323781.57
$ws.Range("T2").Value = 114323781.57

$ws.Range("V2").Value = 13852909.83
$ws.Range("W2").Value = 18249371.06
$ws.Range("X2").Value = 3515562.64
$ws.Range("Y2").Value = -53722201.8
$ws.Range("Z2").Value = -54081440.17
$ws.Range("AA2").Value = -9956410.66

$ws.Range("AG2").Value = 1084943.82

$ws.Range("AS2").Value = -46503129.51

# These ratio cells become blank (empty text) in the new data.
# Assign the "force text" apostrophe so the cell resolves to an empty
# string (Text) rather than being dropped back to a generic numeric
# blank, then strip the quote-prefix formatting it introduces so no
# extra style is left behind on the cell.
foreach ($addr in @("R2", "U2", "AP2", "AQ2", "AR2", "AT2")) {
    $rng = $ws.Range($addr)
    $rng.Value = "'"
    $rng.ClearFormats()
}
